$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column H (Tempo) and column L (Data) carry a specific number format
# (time / date). Copy the formatting down from an already-formatted row
# so the existing style entries are reused instead of minting new ones.
$ws.Range("H2").Copy($ws.Range("H6"))
$ws.Range("H2").Copy($ws.Range("H7"))
$ws.Range("L2").Copy($ws.Range("L6"))
$ws.Range("L2").Copy($ws.Range("L7"))

# --- Row 6 (A6 already = 5) : fill in the rest of the data row ---
$ws.Range("B6").Value = 24
$ws.Range("C6").Formula = "=3*41309"
$ws.Range("D6").Formula = "=C6*F6"
$ws.Range("E6").Value = 300
$ws.Range("F6").Value = 1.5
$ws.Range("G6").Value = 1459
$ws.Range("H6").Value = 0.042256944444444444
$ws.Range("I6").Value = 6657
$ws.Range("J6").Value = "Vampiro"
$ws.Range("K6").Value = "Normal"
$ws.Range("L6").Value = 46012

# --- Row 7 (A7 already = 6) : fill in the rest of the data row ---
$ws.Range("B7").Value = 20
$ws.Range("C7").Formula = "=3*26209"
$ws.Range("D7").Formula = "=C7*F7"
$ws.Range("E7").Value = 293
$ws.Range("F7").Value = 1.5
$ws.Range("G7").Value = 1311
$ws.Range("H7").Value = 0.03400462962962963
$ws.Range("I7").Value = 6695
$ws.Range("J7").Value = "Vampiro"
$ws.Range("K7").Value = "Normal"
$ws.Range("L7").Value = 46013

# --- Rows 8-27: clear the placeholder "Vampiro" tag from column J ---
for ($r = 8; $r -le 27; $r++) {
    $ws.Range("J$r").Value = ""
}

# --- Update the selected cell shown when the sheet is opened ---
$ws.Range("A8").Select()
